$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for Wins/Losses/Ties, matching style of existing header row (A1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (A1) to the new header cells
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Restore the text values (paste formats should not touch values, but just in case)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

for ($row = 2; $row -le 41; $row++) {
    $ws.Cells.Item($row, 30).Value = 81  # AD
    $ws.Cells.Item($row, 31).Value = 81  # AE
    $ws.Cells.Item($row, 32).Value = 0   # AF
}
